# Shift all data on Sheet1 one column to the right (A:K -> B:L)
# by inserting a new, blank column before column A, then move the
# selection/active cell to E8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()

$ws.Range("E8").Select()
